# Consolidate each "Slide N (...)" title (and a couple of caption
# textboxes) that were previously split across several <a:r> runs into
# a single run. PowerPoint's TextRange.Text setter is a no-op when the
# new concatenated text equals the text already there (even though the
# existing text lives in multiple runs), so every real change below is
# done in two steps: set a scratch value first (forcing the run list to
# collapse to one run), then set the final value.

$p = $ppt.ActivePresentation

function Set-ShapeText($slideIndex, $shapeIndex, $finalText) {
    $shp = $p.Slides.Item($slideIndex).Shapes.Item($shapeIndex)
    $tr = $shp.TextFrame.TextRange
    $tr.Text = "__tmp__"
    $tr.Text = $finalText
}

# Titles (shape 1 on every slide)
Set-ShapeText 1  1 "Slide 1 (Content)"
Set-ShapeText 2  1 "Slide 2 (Content)"
Set-ShapeText 3  1 "Slide 3 (Content)"
Set-ShapeText 4  1 "Slide 4 (Content)"
Set-ShapeText 5  1 "Slide 5 (Two Content)"
Set-ShapeText 6  1 "Slide 6 (Two Content Right)"
Set-ShapeText 7  1 "Slide 7 (Content with Caption)"
Set-ShapeText 8  1 "Slide 8 (Comparison)"
Set-ShapeText 9  1 "Slide 9 (Content)"
Set-ShapeText 10 1 "Slide 10 (Content)"
Set-ShapeText 11 1 "Slide 11 (Content)"
Set-ShapeText 12 1 "Slide 12 (Content)"

# "an image" / "An image" caption textboxes
Set-ShapeText 6 3 "an image"
Set-ShapeText 7 4 "An image"
Set-ShapeText 8 4 "An image"
